$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 166-167, pushing the existing historical data
# (previously rows 166-183) down to rows 168-185. This mirrors the
# weekly refresh: a new week's Primera/Segunda quality pair is added at
# the top of this product's price history block.
$ws.Range("A166:R167").Insert()

# Row 166 - new week, Primera quality
$ws.Cells.Item(166, 1).Value = 8
$ws.Cells.Item(166, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(166, 3).Value = "Coquimbo"
$ws.Cells.Item(166, 4).Value = 44504
$ws.Cells.Item(166, 5).Value = 4
$ws.Cells.Item(166, 6).Value = 100114014
$ws.Cells.Item(166, 7).Value = "Betarraga"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 3000
$ws.Cells.Item(166, 11).Value = 450
$ws.Cells.Item(166, 12).Value = 500
$ws.Cells.Item(166, 13).Value = 475
$ws.Cells.Item(166, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(166, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(166, 16).Value = 158
$ws.Cells.Item(166, 17).Value = 3
$ws.Cells.Item(166, 18).Value = "Hortaliza"

# Row 167 - new week, Segunda quality
$ws.Cells.Item(167, 1).Value = 8
$ws.Cells.Item(167, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(167, 3).Value = "Coquimbo"
$ws.Cells.Item(167, 4).Value = 44504
$ws.Cells.Item(167, 5).Value = 4
$ws.Cells.Item(167, 6).Value = 100114014
$ws.Cells.Item(167, 7).Value = "Betarraga"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Segunda"
$ws.Cells.Item(167, 10).Value = 1460
$ws.Cells.Item(167, 11).Value = 350
$ws.Cells.Item(167, 12).Value = 400
$ws.Cells.Item(167, 13).Value = 375
$ws.Cells.Item(167, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(167, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(167, 16).Value = 125
$ws.Cells.Item(167, 17).Value = 3
$ws.Cells.Item(167, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
